# Generate Report for Handoff
# The "cf41b5ba-f47c-4eee-bcfc-6d3cbfddb78f.md" file row moves from
# "In Translation" to "Ready for handoff", and a fresh handoff is recorded
# (new "Latest Handoff Datetime" timestamps) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("D3").Value = "2016-03-09 13:53:33"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("D3").Value = "2016-03-09 13:53:36"
